# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Change cell B11 on sheet "Rules" from the text "R40" to the literal text "1".
#
# A plain `.Value = "1"` assignment would let Excel's type-inference treat the
# numeric-looking string as a number (losing the text type), and prefixing it
# with an apostrophe to force text flips on the cell's "quote prefix" flag,
# which allocates a brand new cell style. Instead we stage the literal text in
# a scratch cell via a formula that evaluates to the text string "1", copy it,
# and paste-special just the values into B11 - this keeps the cell's existing
# style/number format untouched while still landing a true text value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$scratch = $ws.Range("Z1")
$scratch.Formula = "=""1"""
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)   # xlPasteValues
$scratch.Clear()
